$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 23:52"

# Update Cataluña row (row 5) values
$ws.Range("B5").Value = 34726
$ws.Range("C5").Value = 15602
$ws.Range("D5").Value = 15586
$ws.Range("E5").Value = 3538
